$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-23 05:48:37"
$ws.Range("O2").Value = "3.5 °C"
$ws.Range("E3").Value = "2026-02-23 05:48:39"
$ws.Range("H3").Value = "42%"
$ws.Range("E4").Value = "2026-02-23 05:48:41"
$ws.Range("H4").Value = "90%"
$ws.Range("O4").Value = "5.5 °C"
$ws.Range("E5").Value = "2026-02-23 05:48:43"
$ws.Range("N5").Value = "1.6 °C 5:06 TU"
$ws.Range("O5").Value = "2.5 °C"
$ws.Range("E6").Value = "2026-02-23 05:48:46"
$ws.Range("J6").Value = "1025.4 hPa"
$ws.Range("N6").Value = "7.7 °C 5:03 TU"
$ws.Range("O6").Value = "9.3 °C"
$ws.Range("E7").Value = "2026-02-23 05:48:49"
$ws.Range("E8").Value = "2026-02-23 05:48:52"
$ws.Range("J8").Value = "1024.8 hPa"
$ws.Range("N8").Value = "11.6 °C 5:22 TU"
$ws.Range("O8").Value = "13.0 °C"
$ws.Range("E9").Value = "2026-02-23 05:48:54"
$ws.Range("E10").Value = "2026-02-23 05:48:56"
$ws.Range("N10").Value = "3.3 °C 5:11 TU"
$ws.Range("O10").Value = "4.2 °C"
$ws.Range("E11").Value = "2026-02-23 05:48:57"
$ws.Range("N11").Value = "1.8 °C 5:04 TU"
$ws.Range("O11").Value = "2.8 °C"
$ws.Range("E12").Value = "2026-02-23 05:48:59"
$ws.Range("O12").Value = "5.5 °C"
$ws.Range("E13").Value = "2026-02-23 05:49:00"
$ws.Range("J13").Value = "1032.2 hPa"
$ws.Range("N13").Value = "-2.7 °C 5:20 TU"
$ws.Range("O13").Value = "-1.1 °C"
$ws.Range("E14").Value = "2026-02-23 05:49:01"
$ws.Range("L14").Value = "21.6 km/h - 311º 5:26 TU"
$ws.Range("E15").Value = "2026-02-23 05:49:03"
$ws.Range("E16").Value = "2026-02-23 05:49:06"
$ws.Range("H16").Value = "19%"
$ws.Range("E17").Value = "2026-02-23 05:49:09"
$ws.Range("H17").Value = "45%"
$ws.Range("N17").Value = "6.3 °C 5:17 TU"
$ws.Range("O17").Value = "7.0 °C"
$ws.Range("E18").Value = "2026-02-23 05:49:11"
$ws.Range("J18").Value = "1026.0 hPa"
$ws.Range("N18").Value = "1.7 °C 5:20 TU"
$ws.Range("O18").Value = "3.0 °C"
$ws.Range("E19").Value = "2026-02-23 05:49:14"
$ws.Range("E20").Value = "2026-02-23 05:49:16"
$ws.Range("H20").Value = "37%"
$ws.Range("N20").Value = "0.1 °C 5:19 TU"
$ws.Range("O20").Value = "3.0 °C"
$ws.Range("E21").Value = "2026-02-23 05:49:18"
$ws.Range("J21").Value = "1029.3 hPa"
$ws.Range("N21").Value = "1.5 °C 5:17 TU"
$ws.Range("O21").Value = "3.6 °C"
$ws.Range("E22").Value = "2026-02-23 05:49:21"
$ws.Range("H22").Value = "26%"
$ws.Range("L22").Value = "21.6 km/h - 348º 5:29 TU"
$ws.Range("E23").Value = "2026-02-23 05:49:24"
$ws.Range("H23").Value = "27%"
$ws.Range("L23").Value = "35.3 km/h - 342º 5:26 TU"
$ws.Range("E24").Value = "2026-02-23 05:49:27"
$ws.Range("N24").Value = "0.5 °C 5:28 TU"
$ws.Range("O24").Value = "2.5 °C"
$ws.Range("E25").Value = "2026-02-23 05:49:29"
$ws.Range("H25").Value = "30%"
$ws.Range("L25").Value = "18.4 km/h - 345º 5:10 TU"
$ws.Range("E26").Value = "2026-02-23 05:49:31"
$ws.Range("J26").Value = "1026.6 hPa"
$ws.Range("E27").Value = "2026-02-23 05:49:34"
$ws.Range("H27").Value = "34%"
$ws.Range("E28").Value = "2026-02-23 05:49:37"
$ws.Range("H28").Value = "93%"
$ws.Range("N28").Value = "2.1 °C 5:09 TU"
$ws.Range("O28").Value = "3.7 °C"
$ws.Range("E29").Value = "2026-02-23 05:49:40"
$ws.Range("N29").Value = "3.2 °C 5:09 TU"
$ws.Range("O29").Value = "4.2 °C"
$ws.Range("E30").Value = "2026-02-23 05:49:42"
$ws.Range("E31").Value = "2026-02-23 05:49:45"
$ws.Range("H31").Value = "50%"
$ws.Range("J31").Value = "1024.4 hPa"
$ws.Range("L31").Value = "56.9 km/h - 335º 5:03 TU"
$ws.Range("E32").Value = "2026-02-23 05:49:48"
$ws.Range("H32").Value = "95%"
$ws.Range("N32").Value = "-1.2 °C 5:24 TU"
$ws.Range("O32").Value = "1.4 °C"
$ws.Range("E33").Value = "2026-02-23 05:49:50"
$ws.Range("N33").Value = "0.7 °C 5:29 TU"
$ws.Range("O33").Value = "2.5 °C"
$ws.Range("E34").Value = "2026-02-23 05:49:53"
$ws.Range("E35").Value = "2026-02-23 05:49:56"
$ws.Range("H35").Value = "43%"
$ws.Range("N35").Value = "7.6 °C 5:20 TU"
$ws.Range("O35").Value = "10.2 °C"
$ws.Range("E36").Value = "2026-02-23 05:49:58"
$ws.Range("E37").Value = "2026-02-23 05:50:01"
$ws.Range("H37").Value = "79%"
$ws.Range("J37").Value = "1029.5 hPa"
$ws.Range("O37").Value = "3.9 °C"
$ws.Range("E38").Value = "2026-02-23 05:50:03"
$ws.Range("H38").Value = "72%"
$ws.Range("O38").Value = "6.4 °C"
$ws.Range("E39").Value = "2026-02-23 05:50:06"
$ws.Range("O39").Value = "3.5 °C"
$ws.Range("E40").Value = "2026-02-23 05:50:08"
$ws.Range("J40").Value = "1029.7 hPa"
$ws.Range("N40").Value = "0.5 °C 5:29 TU"
$ws.Range("O40").Value = "2.0 °C"
$ws.Range("E41").Value = "2026-02-23 05:50:11"
$ws.Range("J41").Value = "1025.0 hPa"
$ws.Range("K41").Value = "-0.1 MJ/m2"
$ws.Range("E42").Value = "2026-02-23 05:50:13"
$ws.Range("E43").Value = "2026-02-23 05:50:15"
$ws.Range("N43").Value = "2.1 °C 5:29 TU"
$ws.Range("O43").Value = "3.9 °C"
$ws.Range("E44").Value = "2026-02-23 05:50:18"
$ws.Range("E45").Value = "2026-02-23 05:50:21"
$ws.Range("E46").Value = "2026-02-23 05:50:23"
$ws.Range("N46").Value = "0.8 °C 5:01 TU"
$ws.Range("O46").Value = "2.1 °C"
